# Update "paises.xlsx" worksheet: refresh the "last updated" timestamp
# and refresh the per-country case counters (COVID-19 country / Spain
# provincia stats snapshot) per the commit "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "Datos actualizados..." timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 00:22"

# --- Update country statistics rows (B:H = Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 558447
$ws.Range("C4").Value = 25568
$ws.Range("D4").Value = 31986
$ws.Range("E4").Value = 504470
$ws.Range("G4").Value = 1414
$ws.Range("H4").Value = 21991

# Row 5: Espana
$ws.Range("B5").Value = 166831
$ws.Range("C5").Value = 3804
$ws.Range("E5").Value = 87231
$ws.Range("G5").Value = 603
$ws.Range("H5").Value = 17209

# Row 8: Alemania
$ws.Range("B8").Value = 127854
$ws.Range("C8").Value = 2402
$ws.Range("E8").Value = 64532
$ws.Range("G8").Value = 151
$ws.Range("H8").Value = 3022

# Row 10: China
$ws.Range("C10").Value = 0

# Row 17: Brasil
$ws.Range("B17").Value = 22192
$ws.Range("C17").Value = 1230
$ws.Range("E17").Value = 20796

# Row 50: Colombia
$ws.Range("B50").Value = 2776
$ws.Range("C50").Value = 67
$ws.Range("D50").Value = 270
$ws.Range("E50").Value = 2397
$ws.Range("G50").Value = 9
$ws.Range("H50").Value = 109

# Row 165: Libia
$ws.Range("D165").Value = 9
$ws.Range("E165").Value = 15

# Row 169: Maldivas
$ws.Range("D169").Value = 14
$ws.Range("E169").Value = 6
